$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new trailing columns: CheckinDate (F) and CheckoutDate (G) ---

# Header row (row 1): new header text, styled to match the existing bold /
# yellow-fill header cells (A1:E1).
$ws.Range("F1").Value2 = "CheckinDate"
$ws.Range("G1").Value2 = "CheckoutDate"

$headerRange = $ws.Range("F1:G1")
$headerRange.Interior.Color = 65535   # same yellow fill as the existing header
$headerRange.Font.Bold = $true        # same bold font as the existing header

# Data row 2: check-in / check-out dates (stored as Excel serial date
# numbers, matching how the rest of the sheet stores plain numeric values),
# formatted as short dates.
$ws.Range("F2").Value2 = 45628
$ws.Range("G2").Value2 = 45537
$ws.Range("F2:G2").NumberFormat = "m/d/yy"

# Column widths for the two new columns.
$ws.Columns.Item(6).ColumnWidth = 12.92
$ws.Columns.Item(7).ColumnWidth = 13.42

# Row 3 (Sachin / Ingale / ...) is left untouched - only columns F and G
# were added to the sheet, nothing in the existing A:E data changed.

# Update the active selection to J9, matching the saved view state.
[void]$ws.Range("J9").Select()
